$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Move "Bogen 11" (shape id 12) slightly before grouping so the new
# group's bounding box (chOff/chExt) comes out matching the target.
$bogen11 = $s.Shapes.Item("Bogen 11")
$bogen11.Left = 176.1617
$bogen11.Top = 338.3904

# Collect every shape currently on the slide and group them together.
$names = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $names += $s.Shapes.Item($i).Name
}
$range = $s.Shapes.Range($names)
$grp = $range.Group()
$grp.Name = "Gruppieren 45"
